# Auto-generated edit script: updates columns E (depth_mean) and I (logh)
# across multiple worksheets per the commit diff ("rolling down the river").
$wb = $excel.ActiveWorkbook

# --- Sheet "13" (18 cell updates) ---
$ws = $wb.Worksheets.Item("13")
$ws.Range("E4").Value = 0.556155956335123
$ws.Range("I4").Value = -0.254803407197853
$ws.Range("E7").Value = 0.0569053265008198
$ws.Range("I7").Value = -1.24484708050043
$ws.Range("E8").Value = 0.326470276576917
$ws.Range("I8").Value = -0.486156352846017
$ws.Range("E9").Value = 0.0868903014673715
$ws.Range("I9").Value = -1.06102869598244
$ws.Range("E10").Value = 0.0689949376511942
$ws.Range("I10").Value = -1.16118277347733
$ws.Range("E11").Value = 0.141394155033208
$ws.Range("I11").Value = -0.849568543080074
$ws.Range("E12").Value = 0.0770541272743987
$ws.Range("I12").Value = -1.11320409407092
$ws.Range("E13").Value = 0.188932431571157
$ws.Range("I13").Value = -0.723693485998537
$ws.Range("E14").Value = 0.0388449810466055
$ws.Range("I14").Value = -1.41066508611846

# --- Sheet "15" (14 cell updates) ---
$ws = $wb.Worksheets.Item("15")
$ws.Range("E2").Value = 0.33357194499513
$ws.Range("I2").Value = -0.476810482753012
$ws.Range("E3").Value = 0.103246752096473
$ws.Range("I3").Value = -0.986123601354522
$ws.Range("E4").Value = 0.208368019852923
$ws.Range("I4").Value = -0.681168935406929
$ws.Range("E5").Value = 0.115660292693498
$ws.Range("I5").Value = -0.936815712992044
$ws.Range("E7").Value = 0.118618295568362
$ws.Range("I7").Value = -0.925848320654844
$ws.Range("E8").Value = 0.32669127366084
$ws.Range("I8").Value = -0.485862465927253
$ws.Range("E9").Value = 0.148399542142312
$ws.Range("I9").Value = -0.828567438985388

# --- Sheet "5" (16 cell updates) ---
$ws = $wb.Worksheets.Item("5")
$ws.Range("E3").Value = 0.0314608700934574
$ws.Range("I3").Value = -1.50222927053785
$ws.Range("E4").Value = 0.298854491442191
$ws.Range("I4").Value = -0.524540212830461
$ws.Range("E5").Value = 0.071022433807391
$ws.Range("I5").Value = -1.1486044493111
$ws.Range("E7").Value = 0.00408688713481306
$ws.Range("I7").Value = -2.38860735578112
$ws.Range("E8").Value = 0.292886326762752
$ws.Range("I8").Value = -0.533300902644435
$ws.Range("E10").Value = 0.100905322589866
$ws.Range("I10").Value = -0.996085924838993
$ws.Range("E11").Value = 0.292268158933349
$ws.Range("I11").Value = -0.534218496147122
$ws.Range("E12").Value = 0.0398067631956158
$ws.Range("I12").Value = -1.40004313473525

# --- Sheet "5a" (2 cell updates) ---
$ws = $wb.Worksheets.Item("5a")
$ws.Range("E8").Value = 0.578629035198701
$ws.Range("I8").Value = -0.237599777549088

# --- Sheet "6" (15 cell updates) ---
$ws = $wb.Worksheets.Item("6")
$ws.Range("E3").Value = 0.456824240921209
$ws.Range("I3").Value = -0.340250858748654
$ws.Range("E4").Value = 0.38887615975603
$ws.Range("I4").Value = -0.410188680672972
$ws.Range("E6").Value = 0.810515482511753
$ws.Range("I6").Value = -0.0912386848189221
$ws.Range("E8").Value = -0.0509851313580177
$ws.Range("E9").Value = 0.0410550448373626
$ws.Range("I9").Value = -1.38663346922879
$ws.Range("E11").Value = 0.190986898896875
$ws.Range("I11").Value = -0.718996422971846
$ws.Range("E12").Value = 0.0829762877707573
$ws.Range("I12").Value = -1.08104599872781
$ws.Range("E14").Value = 0.186952867332309
$ws.Range("I14").Value = -0.728267869600161

# --- Sheet "6a" (8 cell updates) ---
$ws = $wb.Worksheets.Item("6a")
$ws.Range("E4").Value = 0.226952277503813
$ws.Range("I4").Value = -0.64406545468198
$ws.Range("E5").Value = 0.258582477591671
$ws.Range("I5").Value = -0.587400907696447
$ws.Range("E6").Value = 0.185682797933314
$ws.Range("I6").Value = -0.731228328401752
$ws.Range("E8").Value = 0.346783836740215
$ws.Range("I8").Value = -0.459941152770569

# --- Sheet "7" (12 cell updates) ---
$ws = $wb.Worksheets.Item("7")
$ws.Range("E4").Value = 0.23773763726534
$ws.Range("I4").Value = -0.623902057806076
$ws.Range("E5").Value = 0.21568387517109
$ws.Range("I5").Value = -0.666182322148197
$ws.Range("E6").Value = 0.184310511172043
$ws.Range("I6").Value = -0.734449896393269
$ws.Range("E8").Value = 0.398938311717069
$ws.Range("I8").Value = -0.399094254569193
$ws.Range("E9").Value = 0.110887957954024
$ws.Range("I9").Value = -0.955115614158115
$ws.Range("E10").Value = 0.338755228520558
$ws.Range("I10").Value = -0.470113992943802

# --- Sheet "9" (13 cell updates) ---
$ws = $wb.Worksheets.Item("9")
$ws.Range("E4").Value = -0.0718485345186965
$ws.Range("E6").Value = -0.14322728968482
$ws.Range("E7").Value = 0.0430076588941056
$ws.Range("I7").Value = -1.36645419746087
$ws.Range("E8").Value = -0.089021967023526
$ws.Range("E9").Value = 0.0322295008782968
$ws.Range("I9").Value = -1.49174642013085
$ws.Range("E10").Value = 0.146809310279103
$ws.Range("I10").Value = -0.833246401677268
$ws.Range("E11").Value = 0.1019830414486
$ws.Range("I11").Value = -0.991472040175493
$ws.Range("E13").Value = 0.262140092531933
$ws.Range("I13").Value = -0.581466551606639
